$wb = $excel.ActiveWorkbook

# The workbook starts with zero worksheets; add the single sheet described
# by the target (name "Sheet0") and populate A1 with the shared-string value.
$ws = $wb.Worksheets.Add()
$ws.Name = "Sheet0"
$ws.Range("A1").Value = "Mi primera celda con valor"
